# Apply the DaySale 2025-09-17 update:
#  1. DEXAMETHASONE-AMRIYA 8MG/2ML 3 AMP. row: stock 4:2 -> 4:1, sale price 0.0000 -> 11.8800, txns 0:0 -> 0:1
#  2. New item row "VOLTAREN 75MG/3ML 3 AMP." inserted before "VOLTAREN 75MG/3ML 6 AMP."
#  3. New item row "زيت حرير 100مل" inserted before "سرنجات 3 سم"
#  4. سرنجات 5 سم row: sale price 9.0000 -> 12.0000, txns 3:0 -> 4:0
#  5. Grand total updated, and timestamp footer updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Insert-ItemRow {
    param([int]$RowIndex, [string]$Name, [string]$Stock, [string]$OrderLimit, [string]$Price, [string]$SalePrice, [string]$Txns, [double]$RowHt)

    # Push everything at/after $RowIndex down by one row.
    $ws.Rows("$($RowIndex):$($RowIndex)").Insert()

    # Copy formatting (fonts, fills, borders, number formats) from the row that
    # just got pushed down (RowIndex+1) into the freshly inserted blank row.
    $srcRow = $RowIndex + 1
    $ws.Rows("$($srcRow):$($srcRow)").Copy()
    $ws.Rows("$($RowIndex):$($RowIndex)").PasteSpecial($xlPasteFormats)
    $ws.Application.CutCopyMode = $false
    $ws.Rows("$($RowIndex):$($RowIndex)").RowHeight = $RowHt

    # Recreate the merged regions used by every item row.
    $ws.Range("A$($RowIndex):B$($RowIndex)").Merge()
    $ws.Range("C$($RowIndex):G$($RowIndex)").Merge()
    $ws.Range("H$($RowIndex):K$($RowIndex)").Merge()
    $ws.Range("L$($RowIndex):M$($RowIndex)").Merge()
    $ws.Range("N$($RowIndex):O$($RowIndex)").Merge()

    $ws.Range("C$RowIndex").Value = $Name
    $ws.Range("H$RowIndex").Value = $Stock
    $ws.Range("L$RowIndex").Value = $OrderLimit
    $ws.Range("N$RowIndex").Value = $Price
    $ws.Range("P$RowIndex").Value = $SalePrice
    $ws.Range("Q$RowIndex").Value = $Txns
}

# --- 1. DEXAMETHASONE-AMRIYA row (row 19, unaffected by later inserts) ---
$ws.Range("H19").Value = "4:1"
$ws.Range("P19").Value = "11.8800"
$ws.Range("Q19").Value = "0:1"

# --- 2. Insert "VOLTAREN 75MG/3ML 3 AMP." before row 54 ("VOLTAREN 75MG/3ML 6 AMP.") ---
Insert-ItemRow 54 "VOLTAREN 75MG/3ML 3 AMP." "4:0" "1" "51.00" "16.8300" "0:1" 25.5

# --- 3. Insert "زيت حرير 100مل" before (what is now, post previous insert) row 62 ("سرنجات 3 سم") ---
Insert-ItemRow 62 "زيت حرير 100مل" "0:0" "0" "50.00" "50.0000" "1:0" 25.5

# --- 4. سرنجات 5 سم row, now at row 64 ---
$ws.Range("P64").Value = "12.0000"
$ws.Range("Q64").Value = "4:0"

# --- Renumber the sequential "م" index column (A7:A71) now that two rows were added ---
for ($r = 7; $r -le 71; $r++) {
    $ws.Range("A$r").Value = $r - 6
}

# --- Grand total row moved from 70 to 72 ---
$ws.Range("P72").Value = 3667.24

# --- Footer timestamp row moved from 71 to 73 ---
$ws.Range("A73").Value = "Wednesday, 17 September, 2025 7:21 PM"
